# Tasks.xlsx update
# - CM4 task text is replaced with a question for G, status set to Pending
# - A new CM5 row is inserted carrying over CM4's old task text
# - CI5 task text is updated to describe connecting the Inspection Comment button
# - The "CC" section header text is renamed from "Copy Commentary" to "Copy Comment"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after CM4 (row 152), before the CI section header (row 154).
# Excel COM Insert() shifts row 153 and everything below down by one, and carries
# formatting down from the row above - matching rows 154-209 becoming 155-210.
$ws.Rows(153).Insert()

# CM4 (row 152): repurpose as a clarifying question, mark it Pending
$ws.Range("C152").Value = "Ask G to clarify, what launches the Comment Popup?"
$ws.Range("D152").Value = "Pending"

# New CM5 row (row 153) takes over CM4's previous task text
$ws.Range("A153").Value = "CM5"
$ws.Range("C153").Value = "Connect Shell buttons to show the View"

# CI5 (now row 160 after the insert) gets an updated task description
$ws.Range("C160").Value = "Connect Inspection's 'Inspection Comment' button to launch  popup."

# CC section header (now row 162 after the insert) renamed
$ws.Range("C162").Value = "Copy Comment"

# Refresh the view: clear the previous scroll-pin/selection and select the range
# that was left selected in the saved workbook.
$ws.Range("A47:E210").Select()
